# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (column A) used 4 emoji characters as status icons:
#   📕 (red book)    -> "-3"
#   📘 (blue book)   -> "⚠️" (warning sign)
#   📙 (orange book) -> "+3"
#   📗 (green book)  -> "✅" (check mark)
#
# Replace every occurrence of these emoji with the corresponding plain-text /
# alternate-emoji replacement, keeping everything else on the sheet untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$eRed    = "📕"
$eBlue   = "📘"
$eOrange = "📙"
$eGreen  = "📗"

$newRed    = "-3"
$newBlue   = "⚠️"
$newOrange = "+3"
$newGreen  = "✅"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow  = $firstRow + $used.Rows.Count - 1
$firstCol = $used.Column
$lastCol  = $firstCol + $used.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        # NOTE: keep the string literal on the LEFT of -eq. PowerShell
        # coerces the right-hand side to the type of the left-hand side, so
        # "$val -eq $eRed" would wrongly compare equal for boolean / numeric
        # cell values (e.g. TRUE gets coerced from a non-empty string).
        if ($eRed -eq $val) {
            # "-3" looks like a number, so prefix it with an apostrophe to
            # force Excel to keep storing it as text (same as a user typing
            # '-3 into the cell).
            $cell.Value = "'" + $newRed
        } elseif ($eBlue -eq $val) {
            $cell.Value = $newBlue
        } elseif ($eOrange -eq $val) {
            # "+3" also looks numeric, force text the same way as above.
            $cell.Value = "'" + $newOrange
        } elseif ($eGreen -eq $val) {
            $cell.Value = $newGreen
        }
    }
}
